$d = $word.ActiveDocument

# ============================================================
# Paragraph 1 ("Background" -> "The Investor" paragraph)
# ============================================================

# "Background" (bold heading word) -> "The Investor"
$f1 = $d.Range(0, $d.Content.End)
$f1.Find.Execute("Background", $true, $false, $false, $false, $false, $true, 1, $false, "The Investor", 2)

# "Management Consultant" -> "management consultant" (lower-cased job title)
$f2 = $d.Range(0, $d.Content.End)
$f2.Find.Execute("Management Consultant", $true, $false, $false, $false, $false, $true, 1, $false, "management consultant", 2)

# "...if he contributes the limit to his 401k, how long..." -> "...if he invests all of his disposable income, how long..."
$f3 = $d.Range(0, $d.Content.End)
$f3.Find.Execute("if he contributes the limit to his 401k, how long", $true, $false, $false, $false, $false, $true, 1, $false, "if he invests all of his disposable income, how long", 2)

# ============================================================
# Paragraph 2 ("Robert's employer pays for an investment..." paragraph)
# ============================================================

# "investment consultant" -> "investment adviser"
$f4 = $d.Range(0, $d.Content.End)
$f4.Find.Execute("investment consultant", $true, $false, $false, $false, $false, $true, 1, $false, "investment adviser", 2)

# Append the risk-tolerance discussion (previously its own paragraph) onto the end
# of this paragraph, plus the brand new closing sentence about taking maximal risk.
$pEmployer = $d.Paragraphs.Item(6)
$insertPoint = $d.Range($pEmployer.Range.End - 1, $pEmployer.Range.End - 1)
$appendText = " This makes sense to Robert and suggests that he has a competitive advantage over some of the older employees. He has a time horizon that is measure in decades and no immediate cash needs. Therefore, he is in the position to take the maximal amount of risk and reap the rewards in the long term. "
$insertPoint.InsertAfter($appendText)

# ============================================================
# Paragraph 3 (was "This makes sense..."; becomes "Robert asks...")
# ============================================================

$openQuote = [char]8220
$closeQuote = [char]8221

$pThird = $d.Paragraphs.Item(7)
$bm = $d.Bookmarks.Item("_GoBack")

# Replace everything in the paragraph before the _GoBack bookmark.
$beforeBm = $d.Range($pThird.Range.Start, $bm.Start)
$beforeBm.Text = "Robert asks, " + $openQuote + "Which of the 401k options has the highest risk-reward profile?" + $closeQuote + " Jim explains, " + $openQuote + "The Emerging Markets Equity Index has historically had both the best and most volatile performance." + $closeQuote + " His decision made, Robert wants to know how long he should expect to wait bef"

# Replace everything in the paragraph after the (now-reseated) bookmark, up to
# (but excluding) the trailing paragraph mark.
$bm2 = $d.Bookmarks.Item("_GoBack")
$pThirdNow = $d.Paragraphs.Item(7)
$afterBm = $d.Range($bm2.End, $pThirdNow.Range.End - 1)
$afterBm.Text = "ore he reaches his target."

# ============================================================
# Remove the now-empty paragraph and the "Mabel meets..." paragraph
# ============================================================

$pEmpty = $d.Paragraphs.Item(8)
$pEmpty.Range.Delete()

$pMabel = $d.Paragraphs.Item(8)
$pMabel.Range.Delete()
